# "Generate Report for Handback" - refresh the recorded HO/handback
# generation timestamps for the third tracked file
# (7f97a658-14ba-462a-b69c-22ceb76179b0) across the Overview, zh-cn and
# de-de report sheets. These cells are stored as plain text timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 4)
$overview.Range("G4").Value = "2016-10-27 08:03:22"

# zh-cn sheet: "Correspond Handoff Datetime" (H4) / "Correspond Handback DateTime" (K4)
$zhcn.Range("H4").Value = "2016-10-27 08:03:09"
$zhcn.Range("K4").Value = "2016-10-27 08:04:02"

# de-de sheet: "Correspond Handoff Datetime" (H4) / "Correspond Handback DateTime" (K4)
$dede.Range("H4").Value = "2016-10-27 08:03:22"
$dede.Range("K4").Value = "2016-10-27 08:04:20"
